$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns before column D. The existing D:H (incl. their
# formatting/widths) shift right to I:M, matching the "shift" seen in the diff.
$ws.Range("D1:H1").EntireColumn.Insert()

# Set widths for the newly inserted, still-blank D:H columns so the final
# layout is D=29 E=29 F=31 G=29 H=29 (Excel stores ColumnWidth quantized to
# the default font metrics, so 28.14 / 30.14 round-trip to stored 29 / 31).
$ws.Range("D1").EntireColumn.ColumnWidth = 28.14
$ws.Range("E1").EntireColumn.ColumnWidth = 28.14
$ws.Range("F1").EntireColumn.ColumnWidth = 30.14
$ws.Range("G1").EntireColumn.ColumnWidth = 28.14
$ws.Range("H1").EntireColumn.ColumnWidth = 28.14

# Populate the new D:H columns with the added (older) quarterly data;
# the pre-existing values now living in I:M need no change.
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

$ws.Range("D9").Value = "1400-10-09 (5)"
$ws.Range("E9").Value = "1400-11-24 (3)"
$ws.Range("F9").Value = "1401-04-21 (11)"
$ws.Range("G9").Value = "1401-04-29 (2)"
$ws.Range("H9").Value = "1401-09-15 (6)"

$ws.Range("D11").Value = 19812
$ws.Range("E11").Value = 26372
$ws.Range("F11").Value = 13461
$ws.Range("G11").Value = 24457
$ws.Range("H11").Value = 22337

$ws.Range("D12").Value = -6688
$ws.Range("E12").Value = -8353
$ws.Range("F12").Value = -6973
$ws.Range("G12").Value = -8879
$ws.Range("H12").Value = -9629

$ws.Range("D13").Value = 13123
$ws.Range("E13").Value = 18020
$ws.Range("F13").Value = 6488
$ws.Range("G13").Value = 15578
$ws.Range("H13").Value = 12707

$ws.Range("D14").Value = -769
$ws.Range("E14").Value = -761
$ws.Range("F14").Value = -993
$ws.Range("G14").Value = -1394
$ws.Range("H14").Value = -1579

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

$ws.Range("D16").Value = 4102
$ws.Range("E16").Value = -308
$ws.Range("F16").Value = 7273
$ws.Range("G16").Value = 135
$ws.Range("H16").Value = 4523

$ws.Range("D17").Value = 16457
$ws.Range("E17").Value = 16951
$ws.Range("F17").Value = 12767
$ws.Range("G17").Value = 14319
$ws.Range("H17").Value = 15652

$ws.Range("D18").Value = -190
$ws.Range("E18").Value = -97
$ws.Range("F18").Value = -212
$ws.Range("G18").Value = -52
$ws.Range("H18").Value = -177

$ws.Range("D19").Value = 195
$ws.Range("E19").Value = 280
$ws.Range("F19").Value = 313
$ws.Range("G19").Value = 157
$ws.Range("H19").Value = 167

$ws.Range("D20").Value = 16462
$ws.Range("E20").Value = 17134
$ws.Range("F20").Value = 12869
$ws.Range("G20").Value = 14424
$ws.Range("H20").Value = 15642

$ws.Range("D21").Value = -1231
$ws.Range("E21").Value = -663
$ws.Range("F21").Value = -78
$ws.Range("G21").Value = -1137
$ws.Range("H21").Value = -527

$ws.Range("D22").Value = 15231
$ws.Range("E22").Value = 16471
$ws.Range("F22").Value = 12790
$ws.Range("G22").Value = 13287
$ws.Range("H22").Value = 15115

$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

$ws.Range("D24").Value = 15231
$ws.Range("E24").Value = 16471
$ws.Range("F24").Value = 12790
$ws.Range("G24").Value = 13287
$ws.Range("H24").Value = 15115

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("D26").Value = 27207
$ws.Range("E26").Value = 22078
$ws.Range("F26").Value = 24566
$ws.Range("G26").Value = 25718
$ws.Range("H26").Value = 22988

$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
